# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# ---- NewLoanInput: becomes the active / selected sheet, selection -> B2
$wsNewLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsNewLoanInput.Activate()
$wsNewLoanInput.Range("B2").Select()

# ---- Summary: A4/B4 value fix, selection -> A7
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 5.56
$wsSummary.Range("B4").Value = 5.56
$wsSummary.Range("A7").Select()

# ---- Repayment schedule: disbursement-row figures corrected, fee total row
#      corrected, and O2 (stray empty cell) shifted to P2, selection -> K2
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("I2").Value = 5.56
$wsRepay.Range("K2").Value = 5.56
$wsRepay.Range("L2").Value = 5.56
$wsRepay.Range("O2").Copy($wsRepay.Range("P2"))
$wsRepay.Range("O2").Clear()

$wsRepay.Range("B3").Value = 45
$wsRepay.Range("C3").Value = 42050
$wsRepay.Range("F3").Value = 813.75
$wsRepay.Range("G3").Value = 4186.25
$wsRepay.Range("H3").Value = 73.97
$wsRepay.Range("I3").Value = 0

$wsRepay.Range("B4").Value = 14
$wsRepay.Range("F4").Value = 864.71
$wsRepay.Range("H4").Value = 23.01

$wsRepay.Range("K2").Select()

# ---- Transactions: no longer the active tab, values corrected, selection -> F16
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 16
$wsTrans.Range("E2").Value = 5.56
$wsTrans.Range("H2").Value = 5.56
$wsTrans.Range("A3").Value = 15
$wsTrans.Range("F16").Select()

# ---- Edit Repayment Schedule: disbursement date corrected, selection -> E5
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Range("B4").Value = 42050
$wsEdit.Range("E5").Select()

# Re-activate NewLoanInput last so it ends up the active/tabSelected sheet
$wsNewLoanInput.Activate()
$wsNewLoanInput.Range("B2").Select()
